$wb = $excel.ActiveWorkbook

# Hunk 0: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 8967
$ws.Cells.Item(40, 10).Value = 9772.727999999999
$ws.Cells.Item(40, 12).Value = 9772.727999999999
$ws.Cells.Item(40, 14).Value = -10122.728

# Hunk 1: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(103, 8).Value = 3052.9
$ws.Cells.Item(103, 9).Value = 488.5
$ws.Cells.Item(103, 11).Value = 1465.5
$ws.Cells.Item(103, 13).Value = -879.5

# Hunk 2: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1941.7291
$ws.Cells.Item(132, 9).Value = 1545.8572
$ws.Cells.Item(132, 11).Value = 4637.571599999999
$ws.Cells.Item(132, 13).Value = -2107.571599999999

# Hunk 3: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 3664.36
$ws.Cells.Item(137, 10).Value = 3903.3125
$ws.Cells.Item(137, 12).Value = 11709.9375
$ws.Cells.Item(137, 14).Value = -16809.9375

# Hunk 4: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2515.303
$ws.Cells.Item(32, 9).Value = 2802.6924
$ws.Cells.Item(32, 11).Value = 2802.6924
$ws.Cells.Item(32, 13).Value = -2515.6924

# Hunk 5: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 5995.893
$ws.Cells.Item(61, 9).Value = 5391.0415
$ws.Cells.Item(61, 11).Value = 5391.0415
$ws.Cells.Item(61, 13).Value = -5179.0415

# Hunk 6: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 9526580
$ws.Cells.Item(74, 9).Value = 11906312
$ws.Cells.Item(74, 10).Value = 7654.2856
$ws.Cells.Item(74, 11).Value = 11906312
$ws.Cells.Item(74, 12).Value = 7654.2856
$ws.Cells.Item(74, 13).Value = -11905438
$ws.Cells.Item(74, 14).Value = -9402.285599999999

# Hunk 7: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 9526580
$ws.Cells.Item(77, 9).Value = 11906312
$ws.Cells.Item(77, 10).Value = 7654.2856
$ws.Cells.Item(77, 11).Value = 59531560
$ws.Cells.Item(77, 12).Value = 38271.428
$ws.Cells.Item(77, 13).Value = -59527192
$ws.Cells.Item(77, 14).Value = -47007.428

# Hunk 8: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 4549.5
$ws.Cells.Item(122, 9).Value = 3999.0908
$ws.Cells.Item(122, 10).Value = 5222.222
$ws.Cells.Item(122, 11).Value = 11997.2724
$ws.Cells.Item(122, 12).Value = 15666.666
$ws.Cells.Item(122, 13).Value = -9547.2724
$ws.Cells.Item(122, 14).Value = -20566.666

# Hunk 9: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 5978.58
$ws.Cells.Item(132, 9).Value = 4919.317
$ws.Cells.Item(132, 11).Value = 14757.951
$ws.Cells.Item(132, 13).Value = -12227.951

# Hunk 10: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 5995.893
$ws.Cells.Item(136, 9).Value = 5391.0415
$ws.Cells.Item(136, 11).Value = 16173.1245
$ws.Cells.Item(136, 13).Value = -13623.1245

# Hunk 11: sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1734.8276
$ws.Cells.Item(20, 10).Value = 1610.9445
$ws.Cells.Item(20, 12).Value = 1610.9445
$ws.Cells.Item(20, 14).Value = -2104.9445

# Hunk 12: sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2993.0527
$ws.Cells.Item(99, 9).Value = 2804.25
$ws.Cells.Item(99, 11).Value = 2804.25
$ws.Cells.Item(99, 13).Value = -1306.25

# Hunk 13: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 2891.8948
$ws.Cells.Item(22, 9).Value = 2429.8
$ws.Cells.Item(22, 10).Value = 3405.3333
$ws.Cells.Item(22, 11).Value = 2429.8
$ws.Cells.Item(22, 12).Value = 3405.3333
$ws.Cells.Item(22, 13).Value = -2079.8
$ws.Cells.Item(22, 14).Value = -4105.3333

# Hunk 14: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(25, 8).Value = 2803.5
$ws.Cells.Item(25, 10).Value = 2803.5
$ws.Cells.Item(25, 12).Value = 2803.5
$ws.Cells.Item(25, 14).Value = -3151.5

# Hunk 15: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 2502.4736
$ws.Cells.Item(99, 9).Value = 2241.5
$ws.Cells.Item(99, 11).Value = 2241.5
$ws.Cells.Item(99, 13).Value = -743.5

# Hunk 16: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 2502.4736
$ws.Cells.Item(126, 9).Value = 2241.5
$ws.Cells.Item(126, 11).Value = 6724.5
$ws.Cells.Item(126, 13).Value = -4254.5

# Hunk 17: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 257.625
$ws.Cells.Item(26, 9).Value = 362
$ws.Cells.Item(26, 10).Value = 83.666664
$ws.Cells.Item(26, 11).Value = 1086
$ws.Cells.Item(26, 12).Value = 250.999992
$ws.Cells.Item(26, 13).Value = -798
$ws.Cells.Item(26, 14).Value = -826.999992

# Hunk 18: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(133, 8).Value = 17999.666
$ws.Cells.Item(133, 10).Value = 9999
$ws.Cells.Item(133, 12).Value = 29997
$ws.Cells.Item(133, 14).Value = -40117

# Hunk 19: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(134, 8).Value = 1884.4
$ws.Cells.Item(134, 10).Value = 3333
$ws.Cells.Item(134, 12).Value = 9999
$ws.Cells.Item(134, 14).Value = -20139

# Hunk 20: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(136, 8).Value = 2098
$ws.Cells.Item(136, 10).Value = 5555
$ws.Cells.Item(136, 12).Value = 16665
$ws.Cells.Item(136, 14).Value = -26865

# Hunk 21: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 2005412.2
$ws.Cells.Item(137, 9).Value = 10000000
$ws.Cells.Item(137, 11).Value = 30000000
$ws.Cells.Item(137, 13).Value = -29994900

# Hunk 22: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(138, 8).Value = 7369.2856
$ws.Cells.Item(138, 10).Value = 8275.833000000001
$ws.Cells.Item(138, 12).Value = 24827.499
$ws.Cells.Item(138, 14).Value = -35107.499

# Hunk 23: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 6078.273
$ws.Cells.Item(139, 10).Value = 10310.667
$ws.Cells.Item(139, 12).Value = 30932.001
$ws.Cells.Item(139, 14).Value = -41212.001

# Hunk 24: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(141, 8).Value = 12857.25
$ws.Cells.Item(141, 9).Value = 714.5
$ws.Cells.Item(141, 11).Value = 2143.5
$ws.Cells.Item(141, 13).Value = 3036.5

# Hunk 25: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 2112.9333
$ws.Cells.Item(97, 9).Value = 1549.4445
$ws.Cells.Item(97, 10).Value = 2958.1667
$ws.Cells.Item(97, 11).Value = 1549.4445
$ws.Cells.Item(97, 12).Value = 2958.1667
$ws.Cells.Item(97, 13).Value = -1053.4445
$ws.Cells.Item(97, 14).Value = -3950.1667

# Hunk 26: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2757.6667
$ws.Cells.Item(102, 9).Value = 1797.7
$ws.Cells.Item(102, 10).Value = 3957.625
$ws.Cells.Item(102, 11).Value = 1797.7
$ws.Cells.Item(102, 12).Value = 3957.625
$ws.Cells.Item(102, 13).Value = -175.7
$ws.Cells.Item(102, 14).Value = -7201.625

# Hunk 27: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 1202.5518
$ws.Cells.Item(107, 9).Value = 673.9167
$ws.Cells.Item(107, 10).Value = 3740
$ws.Cells.Item(107, 11).Value = 673.9167
$ws.Cells.Item(107, 12).Value = 3740
$ws.Cells.Item(107, 13).Value = 1246.0833
$ws.Cells.Item(107, 14).Value = -7580

# Hunk 28: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4407.2173
$ws.Cells.Item(132, 9).Value = 3567.647
$ws.Cells.Item(132, 11).Value = 10702.941
$ws.Cells.Item(132, 13).Value = -8172.940999999999

# Hunk 29: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3985.3684
$ws.Cells.Item(46, 9).Value = 769.5
$ws.Cells.Item(46, 10).Value = 4842.933
$ws.Cells.Item(46, 11).Value = 769.5
$ws.Cells.Item(46, 12).Value = 4842.933
$ws.Cells.Item(46, 13).Value = -581.5
$ws.Cells.Item(46, 14).Value = -5218.933

# Hunk 30: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 4224.793
$ws.Cells.Item(100, 9).Value = 3869.8948
$ws.Cells.Item(100, 10).Value = 4899.1
$ws.Cells.Item(100, 11).Value = 3869.8948
$ws.Cells.Item(100, 12).Value = 4899.1
$ws.Cells.Item(100, 13).Value = -3328.8948
$ws.Cells.Item(100, 14).Value = -5981.1

# Hunk 31: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4861.357
$ws.Cells.Item(132, 9).Value = 3850.5
$ws.Cells.Item(132, 10).Value = 6680.9
$ws.Cells.Item(132, 11).Value = 11551.5
$ws.Cells.Item(132, 12).Value = 20042.7
$ws.Cells.Item(132, 13).Value = -9021.5
$ws.Cells.Item(132, 14).Value = -25102.7

# Hunk 32: sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 1678
$ws.Cells.Item(96, 9).Value = 1780
$ws.Cells.Item(96, 10).Value = 1514.8
$ws.Cells.Item(96, 11).Value = 1780
$ws.Cells.Item(96, 12).Value = 1514.8
$ws.Cells.Item(96, 13).Value = -407
$ws.Cells.Item(96, 14).Value = -4260.8

# Hunk 33: sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3555.8333
$ws.Cells.Item(132, 9).Value = 2381.7246
$ws.Cells.Item(132, 10).Value = 12557.333
$ws.Cells.Item(132, 11).Value = 7145.1738
$ws.Cells.Item(132, 12).Value = 37671.999
$ws.Cells.Item(132, 13).Value = -4615.1738
$ws.Cells.Item(132, 14).Value = -42731.999
